# Revert "Powerpoint writer: consolidate text run nodes."
#
# For each paragraph whose first run is "The " / "One " followed by a
# second run ("Moon" / "More"), split the leading run so the trailing
# space becomes its own run again:
#   "The "  -> "The" + " "
#   "One "  -> "One" + " "
#
# Re-assigning TextRange.Text on the Characters() sub-range for the word
# (without its trailing space) forces the host to re-split the backing
# <a:r> run at that character boundary, producing a fresh run (with an
# empty <a:rPr/>) for the remaining trailing space - matching the
# pre-consolidation XML shape.

$p = $ppt.ActivePresentation

function Split-TrailingSpaceRun($shape, $word) {
    $tr = $shape.TextFrame.TextRange
    $len = $word.Length
    $chars = $tr.Characters(1, $len)
    $chars.Text = $word
}

# Slide 2, "TextBox 3": "The " + "Moon"
$s2 = $p.Slides.Item(2)
Split-TrailingSpaceRun $s2.Shapes.Item(2) "The"

# Slide 3, "Title 1": "One " + "More"
$s3 = $p.Slides.Item(3)
Split-TrailingSpaceRun $s3.Shapes.Item(1) "One"

# Slide 3, "TextBox 3": "The " + "Moon"
Split-TrailingSpaceRun $s3.Shapes.Item(3) "The"
